$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 202.9
$ws.Range("I5").Value = 243.42857
$ws.Range("J5").Value = 108.333336
$ws.Range("K5").Value = 243.42857
$ws.Range("L5").Value = 108.333336
$ws.Range("M5").Value = -128.42857
$ws.Range("N5").Value = -338.333336
$ws.Range("H28").Value = 1391.4
$ws.Range("J28").Value = 2281.2
$ws.Range("L28").Value = 2281.2
$ws.Range("N28").Value = -3251.2
$ws.Range("H62").Value = 4248.5
$ws.Range("I62").Value = 3997
$ws.Range("K62").Value = 3997
$ws.Range("M62").Value = -3373
$ws.Range("H65").Value = 4248.5
$ws.Range("I65").Value = 3997
$ws.Range("K65").Value = 19985
$ws.Range("M65").Value = -16865
$ws.Range("H100").Value = 2382.9285
$ws.Range("I100").Value = 2591.9092
$ws.Range("K100").Value = 2591.9092
$ws.Range("M100").Value = -2050.9092
$ws.Range("H132").Value = 3497.5
$ws.Range("I132").Value = 3497.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 10492.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -7962.5
$ws.Range("N132").ClearContents()
$ws.Range("H137").Value = 3352.5293
$ws.Range("I137").Value = 2938.0967
$ws.Range("J137").Value = 7635
$ws.Range("K137").Value = 8814.2901
$ws.Range("L137").Value = 22905
$ws.Range("M137").Value = -6264.2901
$ws.Range("N137").Value = -28005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2477.762
$ws.Range("I32").Value = 2294.3171
$ws.Range("K32").Value = 2294.3171
$ws.Range("M32").Value = -2007.3171
$ws.Range("H43").Value = 29962.334
$ws.Range("J43").Value = 29962.334
$ws.Range("L43").Value = 29962.334
$ws.Range("N43").Value = -30588.334
$ws.Range("H61").Value = 3622.96
$ws.Range("I61").Value = 3233.75
$ws.Range("J61").Value = 5179.8
$ws.Range("K61").Value = 3233.75
$ws.Range("L61").Value = 5179.8
$ws.Range("M61").Value = -3021.75
$ws.Range("N61").Value = -5603.8
$ws.Range("H132").Value = 4772.385
$ws.Range("I132").Value = 4707.1665
$ws.Range("K132").Value = 14121.4995
$ws.Range("M132").Value = -11591.4995
$ws.Range("H136").Value = 3622.96
$ws.Range("I136").Value = 3233.75
$ws.Range("J136").Value = 5179.8
$ws.Range("K136").Value = 9701.25
$ws.Range("L136").Value = 15539.4
$ws.Range("M136").Value = -7151.25
$ws.Range("N136").Value = -20639.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 2649.0833
$ws.Range("I64").Value = 1211.6666
$ws.Range("J64").Value = 3128.2222
$ws.Range("K64").Value = 1211.6666
$ws.Range("L64").Value = 3128.2222
$ws.Range("M64").Value = -986.6666
$ws.Range("N64").Value = -3578.2222
$ws.Range("H67").Value = 2649.0833
$ws.Range("I67").Value = 1211.6666
$ws.Range("J67").Value = 3128.2222
$ws.Range("K67").Value = 1211.6666
$ws.Range("L67").Value = 3128.2222
$ws.Range("M67").Value = -431.6666
$ws.Range("N67").Value = -4688.2222
$ws.Range("H80").Value = 2105.7778
$ws.Range("I80").Value = 165.6
$ws.Range("J80").Value = 4531
$ws.Range("K80").Value = 165.6
$ws.Range("L80").Value = 4531
$ws.Range("M80").Value = 832.4
$ws.Range("N80").Value = -6527
$ws.Range("H83").Value = 2105.7778
$ws.Range("I83").Value = 165.6
$ws.Range("J83").Value = 4531
$ws.Range("K83").Value = 828
$ws.Range("L83").Value = 22655
$ws.Range("M83").Value = 4164
$ws.Range("N83").Value = -32639
$ws.Range("H86").Value = 5307.7334
$ws.Range("J86").Value = 11624.75
$ws.Range("L86").Value = 11624.75
$ws.Range("N86").Value = -13870.75
$ws.Range("H89").Value = 5307.7334
$ws.Range("J89").Value = 11624.75
$ws.Range("L89").Value = 58123.75
$ws.Range("N89").Value = -69355.75
$ws.Range("H105").Value = 5827.593
$ws.Range("I105").Value = 5807.1875
$ws.Range("K105").Value = 5807.1875
$ws.Range("M105").Value = -4060.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1238.5
$ws.Range("I16").Value = 1238.5
$ws.Range("K16").Value = 1238.5
$ws.Range("M16").Value = -951.5
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()
$ws.Range("H31").Value = 2629.8572
$ws.Range("I31").Value = 2806.3635
$ws.Range("J31").Value = 1982.6666
$ws.Range("K31").Value = 2806.3635
$ws.Range("L31").Value = 1982.6666
$ws.Range("M31").Value = -2511.3635
$ws.Range("N31").Value = -2572.6666
$ws.Range("H34").Value = 2629.8572
$ws.Range("I34").Value = 2806.3635
$ws.Range("J34").Value = 1982.6666
$ws.Range("K34").Value = 2806.3635
$ws.Range("L34").Value = 1982.6666
$ws.Range("M34").Value = -2604.3635
$ws.Range("N34").Value = -2386.6666
$ws.Range("H38").Value = 3100
$ws.Range("I38").Value = 3100
$ws.Range("K38").Value = 3100
$ws.Range("M38").Value = -2723
$ws.Range("H41").Value = 21659.75
$ws.Range("I41").Value = 13046.333
$ws.Range("J41").Value = 47500
$ws.Range("K41").Value = 13046.333
$ws.Range("L41").Value = 47500
$ws.Range("M41").Value = -12618.333
$ws.Range("N41").Value = -48356
$ws.Range("H46").Value = 3100
$ws.Range("I46").Value = 3100
$ws.Range("K46").Value = 3100
$ws.Range("M46").Value = -2889
$ws.Range("H113").Value = 1238.5
$ws.Range("I113").Value = 1238.5
$ws.Range("K113").Value = 1238.5
$ws.Range("M113").Value = 931.5
$ws.Range("H122").Value = 3755.8823
$ws.Range("I122").Value = 3722.4
$ws.Range("J122").Value = 4007
$ws.Range("K122").Value = 11167.2
$ws.Range("L122").Value = 12021
$ws.Range("M122").Value = -8717.200000000001
$ws.Range("N122").Value = -16921
$ws.Range("H134").Value = 3205.3845
$ws.Range("I134").Value = 3433.6365
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 10300.9095
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -7765.9095
$ws.Range("N134").Value = -10920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1869.3043
$ws.Range("I26").Value = 2000
$ws.Range("K26").Value = 6000
$ws.Range("M26").Value = -5712
$ws.Range("H34").Value = 64003.234
$ws.Range("J34").Value = 70510.336
$ws.Range("L34").Value = 211531.008
$ws.Range("N34").Value = -211699.008
$ws.Range("H98").Value = 974
$ws.Range("J98").Value = 717.5
$ws.Range("L98").Value = 2152.5
$ws.Range("N98").Value = -5148.5
$ws.Range("H114").Value = 1264.8572
$ws.Range("I114").Value = 1119.6
$ws.Range("J114").Value = 1628
$ws.Range("K114").Value = 3358.8
$ws.Range("L114").Value = 4884
$ws.Range("M114").Value = -104.7999999999997
$ws.Range("N114").Value = -11392
$ws.Range("H122").Value = 553.6
$ws.Range("I122").Value = 518.125
$ws.Range("J122").Value = 695.5
$ws.Range("K122").Value = 4663.125
$ws.Range("L122").Value = 6259.5
$ws.Range("M122").Value = -2213.125
$ws.Range("N122").Value = -11159.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1998.5
$ws.Range("I113").Value = 1998.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1998.5
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 171.5
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 3111.2942
$ws.Range("I122").Value = 2965.3572
$ws.Range("J122").Value = 3792.3333
$ws.Range("K122").Value = 8896.071599999999
$ws.Range("L122").Value = 11376.9999
$ws.Range("M122").Value = -6446.071599999999
$ws.Range("N122").Value = -16276.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1134
$ws.Range("I61").Value = 1134
$ws.Range("K61").Value = 1134
$ws.Range("M61").Value = -932
$ws.Range("H100").Value = 1053.625
$ws.Range("I100").Value = 1017.8571
$ws.Range("K100").Value = 1017.8571
$ws.Range("M100").Value = -476.8570999999999
$ws.Range("H113").Value = 1134
$ws.Range("I113").Value = 1134
$ws.Range("K113").Value = 1134
$ws.Range("M113").Value = 1036
$ws.Range("H132").Value = 2002.0769
$ws.Range("I132").Value = 1928.4375
$ws.Range("J132").Value = 2119.9
$ws.Range("K132").Value = 5785.3125
$ws.Range("L132").Value = 6359.700000000001
$ws.Range("M132").Value = -3255.3125
$ws.Range("N132").Value = -11419.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 15943.5
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 15943.5
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 15943.5
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -16709.5
$ws.Range("H85").Value = 15943.5
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 15943.5
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 15943.5
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -18595.5
$ws.Range("H122").Value = 2290.3333
$ws.Range("I122").Value = 2180.3635
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 6541.0905
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -4091.0905
$ws.Range("N122").Value = -15400
$ws.Range("H132").Value = 3240.2424
$ws.Range("I132").Value = 2794.1292
$ws.Range("K132").Value = 8382.3876
$ws.Range("M132").Value = -5852.3876
